$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = '[''en'' ''nl'' ''sv'']'
$ws.Range("C6").Value = '[''es'' ''en'' ''ar'']'
$ws.Range("D6").Value = '[''fr'' ''en'' ''nl'' ''undetected'' ''ar'']'
$ws.Range("E6").Value = '[''ar'' ''fa'']'
$ws.Range("G6").Value = '[''es'' ''undetected'' ''en'' ''so'' ''ca'' ''de'' ''hu'' ''ar'' ''id'']'
$ws.Range("H6").Value = '[''en'' ''da'']'
$ws.Range("I6").Value = '[''en'' ''mk'' ''ar'' ''hr'' ''fr'' ''tr'' ''id'']'
$ws.Range("J6").Value = '[''ar'' ''so'' ''en'' ''fr'']'
$ws.Range("K6").Value = '[''en'' ''ar'' ''undetected'' ''so'' ''cy'']'
$ws.Range("L6").Value = '[''en'' ''ar'' ''undetected'' ''af'']'
$ws.Range("M6").Value = '[''undetected'' ''es'']'
$ws.Range("N6").Value = '[''fr'' ''pt'' ''en'' ''ar'' ''undetected'' ''nl'' ''es'' ''de'' ''uk'' ''ca'' ''it'' ''sw'']'
$ws.Range("O6").Value = '[''es'' ''en'' ''undetected'' ''ca'' ''pt'' ''so'' ''id'' ''de'' ''it'' ''vi'' ''cy'' ''ar'' ''hu''
 ''sw'']'
$ws.Range("P6").Value = '[''en'' ''fr'' ''ar'']'
$ws.Range("Q6").Value = '[''undetected'' ''en'']'
$ws.Range("R6").Value = '[''ar'' ''en'' ''es'' ''hu'']'
$ws.Range("S6").Value = '[''en'' ''pl'' ''id'' ''fr'' ''ar'']'

$ws.Range("B7").Value = '[''Register'' "Let''s" ''🇳🇱'' nan ''Today'' ''Hi'' ''QUARTER-FINAL'' ''Thanks'' ''📢''
 ''Great'' ''Oranje'' ''As'' ''The'' ''It’s'' ''HE'' ''Minister'' ''Are'' ''Tonight,''
 ''Gender-based'' ''Millions'' ''Match'' ''A'' ''Het'' ''Building'' ''During'' ''We''
 ''Today’s'' ''On'' ''Sustainability'']'
$ws.Range("C7").Value = '[''Desde'' nan ''Coincidiendo'' ''Hoy'' ''Conditions'' ''SM'' ''📲'' ''El'' ''✍️🏻'' ''La''
 ''Si'' ''🔸"هنا'' ''Visitors'' ''Celebración'' ''Con'' ''#DiaDelCineEspañol''
 ''📽️Celebramos'' ''Ya'']'
$ws.Range("D7").Value = '[''Je'' nan ''Today'' ''📣'' ''Sadly'' ''All'' ''A'' ''December'' ''Brimming'' ''Our'' ''🧡''
 ''🇶🇦'' ''Foreign'' ''Deputy'' ''De'' ''⚽️'' ''#DEVILTIME'' ''Team'' "It''s" ''Derniers''
 ''11,59'' ''🚨Just'' ''To'' ''👑'' ''#qatarsustainabilityweek'' ''Belgium'' ''The''
 ''Maj.'' ''Getting'' ''❌'' ''Yesterday'' ''Opening'' ''سعادة'' ''Luhansk,'']'
$ws.Range("E7").Value = '[nan ''تتقدم'' ''در'' ''نظر'' ''🔻پیام'' ''رهبر'' ''گوشه'' ''«شکوه'' ''شکوه'' ''و'' ''⚽️''
 ''این'' ''حضور'' ''#تا_پای_جان_برای_ایران'' ''هواداران'' ''🔴'' ''دارندگان''
 ''♦️اطلاعیه/'' ''🔻آنچه'' ''پیرو'' ''وصول'' ''کاروان'' ''ساعاتی'' ''تحدث'' ''گفتگوی'' ''به'']'
$ws.Range("F7").Value = '[nan]'
$ws.Range("G7").Value = '[''3/3'' nan ''#EstrategiaYDiplomaciaPública'' ''#MásOportunidades'' ''¿Sabías''
 ''@LEGORRETA_ARQS'' ''¡Gran'' ''Después'' ''El'' ''Leaving'' ''التهاني'' ''En'' ''🏆'' ''A''
 ''"Mi🇲🇽querido'' ''Vas'' ''🧵1/3'' ''Recibimos'' ''Un'' ''Nos'' ''▶️'' ''Agradecemos''
 ''La'' ''Con'' ''Conoce'' ''📌Si'' ''📢A'' ''📸'' ''#CentroMexicoQatar2022'' ''Agradecimos''
 ''Afición'' ''𝑼𝒏𝒂'' ''📌¡Si'' ''¿Estás'' ''Cerca'' ''Espléndida'' ''📢'' ''Mexican''
 ''Mexicanos/as'' ''✍️'' ''#ComunicadoSRE'' ''📝'' ''#Qatar2022'' ''🎙️@ALFZEGBE:''
 ''Atención'' ''Este'' ''Estrategia'' ''Horas'' ''🇲🇽🗣️'' ''Que'' ''32'' ''Charro''
 ''@Tralexita1'' ''¡Gracias'' ''¡Hoy'' ''¡Entusiasta'' ''📌'' ''⚠️Nuestra'' ''Despite''
 ''▶️Si'' ''⚽️✴️La'' ''Guillermo'' ''Damos'' ''1/2'' ''Para'' ''@RafaelLaveagaR''
 ''Debutamos'' ''🚨Invitamos'' ''Actualización-los'' ''Terminando''
 ''.@miseleccionmx'' ''🤩¡Con'' ''Estás'' ''Compartiendo'' ''Si'' ''1/3'' ''🇲🇽¡𝐀𝐪𝐮𝐢́''
 ''Empezando'' ''سعادة'' ''¡Un'' ''Desde'' ''Seguimos'' ''Colegas'' ''🇦🇷🇲🇽'' ''¿Por''
 ''Recuerda'' ''Y'' ''Toma'' ''✅'' ''¡Bienvenido'' ''Muy'' ''¡Enhorabuena!'' ''H.E''
 ''كلمة'' ''🚨'' ''@AnguloNTN24'' ''¡Evita'' ''2/2🚨'' ''🎉'' ''@miseleccionmx''
 ''PodcastIMR🎙️'' ''Had'' ''¡Se'' ''▶️¿Vienes'' ''¡#LeleMundialista'' ''📢Atención''
 ''❗️Si'' ''@DiplomaciaPubl'' ''➡️'' ''🗣️¡Estamos'' ''¡Los'' ''Hoy'' ''Sofia'' ''De''
 ''#CentroMéxicoQatar2022'' ''¡En'' ''@FIFAWorldCup'' ''Coming'' ''¡Comienza''
 ''¡Felicidades!'' ''¡Felicidades'' ''🎉Entusiasta'' ''Aquí'' ''Es'' ''Welcome''
 ''¡Estamos'' ''🧵Muchas'' ''México'' ''¡Así'' ''¡Faltan'' ''Su'' ''▶️¡Faltan'' ''Mi''
 ''Los'' ''#EnVivo'' ''Inicia'' ''💀🕯️'' ''🧵A'' ''¡Aquí'' ''Ofrecí'' ''¡La'' ''Asistí'' ''📣''
 ''Importantes'' ''🇲🇽▶️'' ''Estamos'' ''Anunciamos'' ''🔶'' ''Gracias,''
 ''#DatoMundialista'' ''“La'' ''¡Falta'' ''Buenos'' ''On'' ''#Mexico'' ''#SeeYouIn2026''
 ''🇲🇽Si'' ''Mexico,'' ''#MiConsulado'' ''Mural'' ''تدشين'' ''🔴سعادة'' ''يسرّنا''
 ''Special'' ''We'' ''💀🌻'' ''من'' ''📣Join'' ''¿Vienes'' ''Yon'' ''⚠️Si'' ''▶️¡Estamos''
 ''🔵"Volemos'' ''🇲🇽🟣'' ''Fue'' ''Las'' ''🇲🇽⚽️Aficionadx'' ''📄'' ''🇲🇽'' ''¿Vendrás''
 ''⚽️¡Ya'' ''Platicamos'' ''🇲🇽¿Vienes'' ''🗣️Faltan'' ''#EnVIVO'' ''🥳Del'' ''Mexicanx''
 ''¿🇲🇽Vives'' ''¿Eres'' ''Me'' ''"في'' ''¡No'' ''🇲🇽Mexicanx'' ''🇲🇽✈️🇶🇦¿Vienes'' ''⌛️'']'
$ws.Range("H7").Value = '[''Very'' ''HE'' ''I'' nan ''Yesterday'' ''#UNDay'' ''⚠️'']'
$ws.Range("I7").Value = '[nan ''Последњи'' ''FIFA'' ''Happy'' ''Embassy'' ''Président'' ''🧡'' ''Амбасада''
 ''Главу'' ''НЕ'' ''🇷🇸'' ''Talking'' ''أهم'' ''Statement'' ''The'' ''Group'' ''On'' ''Данас''
 ''During'' ''#WorldCupQatar2022'' ''Charge'' ''Visitors'' ''У'' ''Many'' ''شهد'' ''In''
 ''One'' ''©️🇷🇸'' ''#Изјава'' ''FM'' ''#Serbia'' "#Bartın''da" ''Shocked'' ''Sun''
 ''Министар'']'
$ws.Range("J7").Value = '[nan ''#المونديال_شرقي'' ''📌'']'
$ws.Range("K7").Value = '[nan ''A'' ''⚽'' ''The'' ''خالص'' ''What'' ''Who'' ''With'' ''This'' ''Enjoyed''
 ''#GreatCreativity'' ''Building'' ''#QNA_Infographic'' ''While'' ''Huge'' ''Our''
 ''England'' ''Busy'' ''Come'' ''Good'' ''Are''
 ''🏴󠁧󠁢󠁥󠁮󠁧󠁿'' ''Thank''
 ''Meet'' ''Citizens'' ''Is'' ''Congratulations'' ''Hats'' ''Today'' ''Interested''
 ''Ahead'' ''Of'' ''رئيس'' ''As'' ''So'' ''Ambassador'' ''Before'' ''PM,'' ''🧵'' ''One''
 ''🇸🇦Saudi'' ''Traveling'' ''You'' "You''ll" ''Just'' ''It’s'' ''Was'' "That''s" ''أمسية''
 ''Matchday'' ''To'' ''من'' ''If'' ''Highlights'' ''اليوم'' ''Don’t'' ''إن'' ''يواجه''
 ''المباراة'' ''صور'' ''From'' ''كجزء'' ''Well'' ''سعدت'' ''Pleased'' ''Great'' ''Pob''
 ''🇺🇸🏴󠁧󠁢󠁷󠁬󠁳󠁿'' ''Best''
 ''Fe'' ''Do'' ''Some'' ''بعد'' ''المنتخب'' ''#LestWeForget'' "It''s" ''كلمة'' ''For''
 ''Less'' ''Residents'' ''Wondering'' ''WATCH:'' ''ظلّ'' ''جلالة'' ''After'' ''UK'' ''كم''
 ''هؤلاء'' ''These'' ''سعادة'' ''My'' ''#StreetChildWorldCup'' ''Planning'' ''It''
 ''الشعار'']'
$ws.Range("L7").Value = '[nan ''Our'' ''Next'' ''Outstanding'' ''Congratulations'' ''لقد'' ''The'' ''Thanks''
 ''Just'' ''#TBT:'' ''And'' ''If'' ''As'' ''It’s'' ''So'' ''Australia'' ''#Aussie!'' ''Let’s''
 ''I'' ''من'' ''Wishing'' ''Thank'' ''أمسية'' ''A'' ''Go'' ''Assistant'' ''It'' ''Australian''
 ''Was'' ''We’re'' ''1️⃣'' ''OPEN'' ''Who'' ''🎧Tune'' ''#FIFAWorldCup'' ''For'' ''Register''
 ''Welcome'' ''Today'' ''What'' ''يسر'' ''مساعد'' "Can''t" ''From'' ''Coming'' ''On''
 ''History'' ''Pleased'' ''EAQC'' ''Only'' ''Qatar'' "I''ll" ''"The'' ''Further'' ''Are''
 ''Great'' ''#Qatar:'' ''Looking'']'
$ws.Range("M7").Value = '[nan]'
$ws.Range("N7").Value = '[nan ''France🇫🇷'' ''La'' ''#العالم_العربي'' ''#Qatar'' ''بمناسبة''
 ''#التنوّع_البيولوجي'' ''Votre'' ''Les'' ''Bravo'' ''Nous'' ''شكرا'' ''MERCI'' ''💔'' ''Au''
 ''Allez,'' ''𝘼𝙇𝙇𝙀𝙕'' ''C’est'' ''مباراة'' ''Match'' ''Voici'' ''حماس'' ''𝙳𝚎𝚛𝚗𝚒𝚎̀𝚛𝚎𝚜''
 ''🇦🇷🇫🇷'' ''#بالفيديو'' ''Et'' ''Un'' ''"On'' ''French'' ''4️⃣'' ''💌'' ''#مجلس_قناة_الكاس''
 ''An'' ''🎥|'' ''🗣'' ''🔴'' ''𝟮𝗲'' ''🏆'' ''الرئيس'' ''مبابي'' ''Direction'' ''𝐎𝐍'' ''Merci''
 ''Encore'' ''Coup'' ''Quelle'' ''انطلاق'' ''Joie'' ''توقعاتكم'' ''Plus'' ''هل'' ''🇫🇷🇶🇦|''
 ''#QNA_Video'' ''Prêts'' ''Deschamps'' ''Une'' ''🫡'' ''Tomorrow,'' ''Demain''
 ''#Audiences'' ''وزيرة'' ''🧨💥'' ''Bonjour'' ''الحضور'' ''Notre'' ''ديشامب'' ''EN'' ''Mais''
 ''⏸️'' ''Cette'' ''ALLEZ'' ''May'' ''À'' ''𝗘𝗡𝗦𝗘𝗠𝗕𝗟𝗘'' ''Échange'' ''Arrivée''
 ''@equipedefrance'' ''🇨🇵أو'' ''Place'' ''🔜'' ''كرم'' ''مبروك'' ''💬'' "C''est" ''السفير''
 ''الف'' ''هدف'' ''Grand'' ''Fact!!!'' ''⚽️'' ''𝑆𝑒𝑢𝑙'' ''𝗢𝗡'' ''𝐔𝐧𝐞'' ''GIMS,'' ''أمسية''
 ''شكرًا'' ''جزيل'' ''Remerciements'' ''Très'' ''Ce'' ''Après'' ''𝘾𝙝𝙖𝙦𝙪𝙚'' "📍C''est"
 ''Toujours'' ''سافرا'' ''Ils'' ''𝙌𝙐𝘼𝙇𝙄𝙁𝙄𝙀𝙎'' ''𝙆𝙔𝙇𝙄𝘼𝘼𝘼𝘼𝘼𝘼𝘼𝘼𝙉'' ''QNA'' ''𝘾𝙤𝙪𝙥''
 ''المشجعون'' ''عودة'' ''زيارة'' ''Déplacement'' ''Retour'' ''🤝'' ''📸|'' ''L’ambiance''
 ''𝙋𝙧𝙚𝙢𝙞𝙚̀𝙧𝙚'' "71''" ''Kyliaaaaaaaan'' ''OLIVIER'' ''Réaction'' "L''ambassade"
 ''Belle'' ''𝐿𝑒'' ''𝗛𝗮̂𝘁𝗲'' ''#Qatar2022'' ''حفل'' ''Cérémonie'' ''وزير'' ''The''
 ''Heureux'' ''لحظة'' ''سعيد'' ''After'' ''Moment'' ''Conférence'' ''1er'' ''👏The''
 ''لقطات'' ''بعد'' ''#BREAKING'' ''وصل'' ''الجماهير'' ''عاجل:'' ''Bienvenue'' ''📸'' ''👋''
 ''#Paris2024'' ''اللقاء'' ''Le'' ''Fier'' ''اختتاماً'' ''مساعد'' ''Vous'' ''قائمة'' ''𝙇𝙖''
 ''Afin'' ''تُعقد'' ''⏳'' ''استضافة'' ''تشرّفت'' ''جان'' ''#سورية'' ''Dans'' ''Son'' ''أود''
 ''#Ukraine'' ''🇪🇺'' ''Two'' ''سعادة'' ''#FRQA50'' ''أسعدتني'' ''مقر'' ''Ravi'' ''To''
 ''Visite'' ''سُعدت'' ''شاركت'' ''"بغداد'' ''22'' ''«'' ''احتفلنا'' ''Célébration''
 ''Plaisir'' ''"متحف'' ''معرض'' ''"سفر"'' ''يكشف'' ''L’expo.'' ''Ukraine'' ''28/08/2024''
 ''H.E'' ''🎥'' ''Striker'' ''🇶🇦'' ''HH'' ''لقاء'' ''سعدني'' ''I'' ''Karim'' ''KARIM''
 ''Delighted'' ''Mission'' ''احتفت'' "L''invitée" ''#UNGA'' ''@QatarTelevision''
 ''أسعدني'' ''🇶🇦📖❤️The'' ''تَمّ'' ''أشكر'' ''@Qatar_Museums'' ''En'' ''#lepetitprince''
 ''Je'' ''قراءة'' ''@LettheWorldread'' ''56'' ''L’exposition'' "☀️Aujourd''hui,"
 ''بالإنابة'' ''زُوروا'' ''أتقدّم'' ''شرّفني'' ''Venez'' ''A'' ''BREAKING'' ''تعرفوا''
 ''بحضور'' ''سيريل'' ''يوري'' ''1/2)'' "L''exposition" ''Sur'' ''On'' ''يلتقي''
 ''Sustainability'']'
$ws.Range("O7").Value = '[''Chau'' nan ''Qatar'' ''🍒🇦🇷🇶🇦'' ''Argentina'' ''Faro'' ''¡Se'' ''@Guille__Nicolas''
 ''Despedimos'' ''Al'' ''Los'' ''🇦🇷'' ''Lionel'' ''Muchas'' ''#Qatar2022''
 ''#NuevaFotoDePerfil'' ''¡PÓNGANLE'' ''Lifting'' ''@Argentina'' ''ADD'' ''¡NO''
 ''OPERATIVO'' ''#Our_Unity_Source_of_Our_Strength'' ''¡Hoy'' ''HAPPY'' ''⚽'' ''¡Que''
 ''Con'' "Couldn''t" ''ℹ️'' ''📢'' ''¡¡@Argentina'' ''🍒'' ''Thank'' ''“Muchaaaachos”''
 ''💙🤍💙'' ''🇶🇦🇦🇷'' ''CELEBRACION'' ''Alojamiento'' ''It'' ''#TodosJuntos'' ''¡#ARG'' ''¡''
 ''¡Así'' ''La'' ''El'' ''¡El'' ''The'' ''🇦🇷🇦🇺'' ''Conditions'' ''Sí,'' ''¡Ya''
 ''Clasificados'' ''“Muchaaaaaachos”'' ''🚈'' ''أهازيج'' ''Nuestro'' ''Bueno,''
 ''¡VAMOS'' ''Publicamos'' ''Guillermo'' ''¡𝐅𝐈𝐄𝐒𝐓𝐀'' ''¡Muchas'' ''Así'' ''Bienvenido''
 ''📸'' ''OFICIAL'' ''¡Gracias'' ''¡Ganó'' ''¡Momento'' ''@ordorica_g'' ''¡¡¡VAMOS''
 ''الجماهير'' ''Your'' ''En'' ''Participamos'' ''DÍA'' ''✈️'' ''Witness'' ''Tras''
 ''Estamos'' ''¿Manijas'' ''🇶🇦'' ''Podés'' ''🤔'' ''🚌'' ''🏬'' ''🇦🇷🇲🇽'' ''#FIFAWorldCup''
 ''¡Llegó'' ''¡Acá'' ''Durante'' ''High'' ''GALA'' ''Faltan'' ''Get'' ''Si''
 ''#SelecciónMayor'' ''Mi'' ''IN'' ''Yendo'' ''📲'' ''¡Bienvenida'' ''¡Ahora'' ''💾''
 ''Sugerimos'' ''IMPORTANTE:'' ''🏆'' ''🎶'' ''Guía'' ''Tickets?'' ''📱'' ''🏟️ESTADIO''
 ''@pepeotegui'' ''This'' ''4️⃣'' ''Un'' ''¿Dijeron'' ''¡La'' ''Chants'' ''Nuestros''
 ''#قنا_إنفوجرافيك'' ''Ya'' ''Este'' ''Bienvenidos'' ''Después'' ''¡Qatar'' ''WATCH''
 ''سعادة'' ''Hoy'' ''Ubicado'' ''Celebramos'' ''#Malvinas'' ''FIESTA'' ''👆🇦🇷👆'' ''Hola''
 ''Around'' ''🫐'' ''A'' ''Para'' ''Desde'' ''Gracias'' ''An'' ''Comparto'' ''Leo'' ''Empezó''
 ''To'' ''Embajador'' ''🇦🇷⚽️'' ''جماهير'' ''Más'' ''Estos'' ''Atardeceres''
 ''Repercusiones'' ''💪🇦🇷'' ''Junto'' ''With'' ''Puerto'' ''Como'' ''🤩'' ''⏰'' ''🍻'' ''Mural''
 ''👇'' ''صحيفة'' ''@LuccoMarcelo'' ''سكالوني'' ''😍'' ''45'' ''¡Felicitaciones'' ''8''
 ''It’s'' ''¡Y'' ''Jueves'' ''#Qatar'' ''🪪'' ''¿Vas'' ''¡FALTA'' ''🍊🍊'' ''@gustrivi'' ''🫐🫐🫐''
 ''Every'' ''¿Todavía'' ''32'' ''¡Fiesta'' ''Located'' ''974'' ''Abu'' ''Those''
 ''Arrangements'' ''Qatari'' ''Fans'' ''#Reconocimiento'' ''¡Si'' ''طرح'' ''Entre''
 ''Cerramos'' ''🇦🇷❄️'' ''Proyección'' ''🇬🇧Two'' ''🗣'' ''Almuerzo'' ''🐮🍊'' ''Is''
 ''Camiseta'' ''Visita'' ''Three'' ''🚨'' ''11'' ''Se'' ''¡Camiseta'' ''🗓'' ''📽'' ''#Maradona''
 ''🔸'' ''Boletín'' ''"Qatar'' ''Y'' ''عاجل'' ''INVITACION:'' ''Brasil'' ''👆😄👆''
 ''¡Felicitaciones!'' ''🐦'' ''Former'' ''Siguen'' ''¡Momentos'' ''خافيير'' ''In''
 ''Tweet'' ''Flag'' ''Aquí'' ''Ceremonia'' ''¡No'' ''One'' ''من'' ''SOLO.'' ''⚽️'' ''Copa''
 ''🔥🇦🇷🔥'' ''Show'' ''بابو'']'
$ws.Range("P7").Value = '[nan ''Today'' ''Congratulations,'' ''Well-done'' ''Listen:'' ''Proud'' ''A'' ''Today,''
 ''In'' ''Il'' ''Engaging'' ''Vous'' ''You'' ''Un'' ''Nous'' ''We'' ''Honourable''
 ''Minister'' ''Le'' ''HOW'' ''Morning'' ''#2022FIFAWorldCup'' ''Consiel'' ''6'' ''The''
 ''26'' ''Canada'' ''13'' ''#QNA_Infographic'' ''Media'' ''Conseil'' ''🇨🇦'' ''20'' ''Notre''
 ''Our'' ''Attention!'' ''30'' ''27'' ''Planning'' ''#FIFACoupeduMonde2022'' ''48'']'
$ws.Range("Q7").Value = '[nan ''Day'']'
$ws.Range("R7").Value = '[nan ''يا'' ''#TBT'' ''The'' ''Our'' ''Today,'' ''What'' ''نتقدم'' ''سيتم'' ''إن'' ''Human''
 ''تهانينا'' ''Congratulations'' ''How'' ''Thanks'' ''بمناسبة'' ''On'' ''حزنّا'' ''يبدأ''
 ''U.S.'' ''Today'' ''نهنئ'' ''American'' ''سيخوض'' ''قال'' ''بالتوفيق'' ''Good''
 ''Friday’s'' ''Pleased'' ''يواجه'' ''المباراة'' ''🇺🇸&amp;🇶🇦'' "We''re" ''Traveling''
 ''Celebrate'' ''Grateful'' ''Exciting'' ''With'' ''GOOOOAAAL!'' ''One'' ''Are'' ''بعد''
 ''Si'' ''If'' "It''s" ''Don’t'' ''El'' ''Qatar’s'' ''In'' ''It'' ''RT'' ''America'' ''Museum''
 ''Calling'' ''#خلف_الكواليس'' ''#BehindTheScenes'' ''أقل'' ''Five'' ''Hello'' ''هذا''
 ''This'' ''يعملون'' ''ساعد'' ''شهر'' ''October'' ''Ready'' ''We’re'' ''Ambassador''
 ''Mental'' ''أقام'' ''Last'' ''Indigenous'' ''من'' ''Many'' ''Qatar'' ''اجتمع'' ''Experts''
 ''متواجدة'' ''Working'']'
$ws.Range("S7").Value = '[''Today,'' nan ''Mikołajki'' "Russia''s" ''Ales'' ''The'' ''#Szczesny💟⚽️'' ''In''
 ''⚠️Paszport'' ''#RobertLewandowski'' ''TO'' ''Za'' ''Oficjalny'' ''⏬️W'' ''⚽️Przed''
 ''Dzień'' ''2️⃣2️⃣'' ''Wwóz'' ''⚽️Wybierasz'' ''⚠️Do'' ''This'' ''During'' ''Z''
 ''⚽️Udostępniliśmy'' ''⚽️'' ''💬'' ''#OTD'' ''#Polki'' ''اللجنة'' ''Modern'' ''لا''
 ''Massive'' "''When" ''How'' "It''s" ''📺PM'' ''On'']'

$ws.Rows.Item(6).EntireRow.AutoFit()
$ws.Rows.Item(7).EntireRow.AutoFit()
